# Generate Report for Handoff
# Updates the localization-status report after a handoff/handback run for
# file "5125f0cd-0be7-4653-9ce4-35f7ce8120a7.md" and records that the
# handback for "3e537be1-bbc8-4711-8dfd-c0e1ac23b1c6.md" (zh-cn) has
# completed (Latest Target File / Latest Handback File now populated).

$wb = $excel.ActiveWorkbook

# --- Overview sheet: bump the "Latest HO Xliff Generate Date" for the
#     5125f0cd-... row (row 7, column G) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G7").Value = "2016-08-28 02:41:33"

# --- zh-cn sheet: handback completed for row 2 (3e537be1-...) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("I2").Value = "3e537be1-bbc8-4711-8dfd-c0e1ac23b1c6.md"
$wsZhCn.Range("J2").Value = "3e537be1-bbc8-4711-8dfd-c0e1ac23b1c6.fff2734d607640bd36765059c09fb28d3bc65cc1.zh-cn.xlf"

# --- zh-cn sheet: new handoff datetime for row 7 (5125f0cd-...) ---
$wsZhCn.Range("H7").Value = "2016-08-28 02:41:28"

# --- de-de sheet: new handoff datetime for row 7 (5125f0cd-...) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H7").Value = "2016-08-28 02:41:33"
